$wb = $excel.ActiveWorkbook

$wsARM = $wb.Worksheets.Item("ARM")
# Row 4
$wsARM.Range("H4").Value = 1115.8334
$wsARM.Range("I4").Value = 1636.125
$wsARM.Range("J4").Value = 855.6875
$wsARM.Range("K4").Value = 1636.125
$wsARM.Range("L4").Value = 855.6875
$wsARM.Range("M4").Value = -1520.125
$wsARM.Range("N4").Value = -1087.6875
# Row 6
$wsARM.Range("H6").Value = 7001
$wsARM.Range("I6").Value = 7001
$wsARM.Range("J6").Value = 0
$wsARM.Range("K6").Value = 7001
$wsARM.Range("L6").Value = 0
$wsARM.Range("M6").Value = -6828
$wsARM.Range("N6").ClearContents()
# Row 9
$wsARM.Range("H9").Value = 0
$wsARM.Range("I9").Value = 0
$wsARM.Range("J9").Value = 0
$wsARM.Range("K9").Value = 0
$wsARM.Range("L9").Value = 0
$wsARM.Range("N9").ClearContents()
# Row 20
$wsARM.Range("H20").Value = 0
$wsARM.Range("I20").Value = 0
$wsARM.Range("J20").Value = 0
$wsARM.Range("K20").Value = 0
$wsARM.Range("L20").Value = 0
$wsARM.Range("N20").ClearContents()
# Row 37
$wsARM.Range("H37").Value = 17719
$wsARM.Range("I37").Value = 0
$wsARM.Range("J37").Value = 17719
$wsARM.Range("K37").Value = 0
$wsARM.Range("L37").Value = 17719
$wsARM.Range("N37").Value = -18265
# Row 44
$wsARM.Range("H44").Value = 0
$wsARM.Range("I44").Value = 0
$wsARM.Range("J44").Value = 0
$wsARM.Range("K44").Value = 0
$wsARM.Range("L44").Value = 0
$wsARM.Range("N44").ClearContents()
# Row 55
$wsARM.Range("H55").Value = 0
$wsARM.Range("I55").Value = 0
$wsARM.Range("J55").Value = 0
$wsARM.Range("K55").Value = 0
$wsARM.Range("L55").Value = 0
$wsARM.Range("N55").ClearContents()
# Row 74
$wsARM.Range("H74").Value = 3849.963
$wsARM.Range("I74").Value = 22000
$wsARM.Range("J74").Value = 2397.96
$wsARM.Range("K74").Value = 22000
$wsARM.Range("L74").Value = 2397.96
$wsARM.Range("M74").Value = -21126
$wsARM.Range("N74").Value = -4145.96
# Row 77
$wsARM.Range("H77").Value = 3849.963
$wsARM.Range("I77").Value = 22000
$wsARM.Range("J77").Value = 2397.96
$wsARM.Range("K77").Value = 110000
$wsARM.Range("L77").Value = 11989.8
$wsARM.Range("M77").Value = -105632
$wsARM.Range("N77").Value = -20725.8
# Row 80
$wsARM.Range("H80").Value = 23600
$wsARM.Range("I80").Value = 0
$wsARM.Range("J80").Value = 23600
$wsARM.Range("K80").Value = 0
$wsARM.Range("L80").Value = 23600
$wsARM.Range("M80").ClearContents()
$wsARM.Range("N80").Value = -25596
# Row 83
$wsARM.Range("H83").Value = 23600
$wsARM.Range("I83").Value = 0
$wsARM.Range("J83").Value = 23600
$wsARM.Range("K83").Value = 0
$wsARM.Range("L83").Value = 70800
$wsARM.Range("M83").ClearContents()
$wsARM.Range("N83").Value = -80784

$wsCRP = $wb.Worksheets.Item("CRP")
# Row 58
$wsCRP.Range("H58").Value = 1294.8611
$wsCRP.Range("I58").Value = 951.1429000000001
$wsCRP.Range("J58").Value = 1513.591
$wsCRP.Range("K58").Value = 951.1429000000001
$wsCRP.Range("L58").Value = 1513.591
$wsCRP.Range("M58").Value = -748.1429000000001
$wsCRP.Range("N58").Value = -1919.591
# Row 105
$wsCRP.Range("H105").Value = 818.5
$wsCRP.Range("I105").Value = 758
$wsCRP.Range("J105").Value = 1000
$wsCRP.Range("K105").Value = 758
$wsCRP.Range("L105").Value = 1000
$wsCRP.Range("M105").Value = 989
$wsCRP.Range("N105").Value = -4494
# Row 132
$wsCRP.Range("H132").Value = 2420.6
$wsCRP.Range("I132").Value = 1973.1428
$wsCRP.Range("J132").Value = 3464.6667
$wsCRP.Range("K132").Value = 5919.428400000001
$wsCRP.Range("L132").Value = 10394.0001
$wsCRP.Range("M132").Value = -3389.428400000001
$wsCRP.Range("N132").Value = -15454.0001
# Row 136
$wsCRP.Range("H136").Value = 1294.8611
$wsCRP.Range("I136").Value = 951.1429000000001
$wsCRP.Range("J136").Value = 1513.591
$wsCRP.Range("K136").Value = 2853.4287
$wsCRP.Range("L136").Value = 4540.772999999999
$wsCRP.Range("M136").Value = -303.4287000000004
$wsCRP.Range("N136").Value = -9640.772999999999
# Row 140
$wsCRP.Range("H140").Value = 19999
$wsCRP.Range("I140").Value = 0
$wsCRP.Range("J140").Value = 19999
$wsCRP.Range("K140").Value = 0
$wsCRP.Range("L140").Value = 19999
$wsCRP.Range("N140").Value = -30359

$wsCUL = $wb.Worksheets.Item("CUL")
# Row 131
$wsCUL.Range("H131").Value = 1627015.4
$wsCUL.Range("I131").Value = 16666886
$wsCUL.Range("J131").Value = 1083.3784
$wsCUL.Range("K131").Value = 50000658
$wsCUL.Range("L131").Value = 3250.1352
$wsCUL.Range("M131").Value = -49995618
$wsCUL.Range("N131").Value = -13330.1352
# Row 132
$wsCUL.Range("H132").Value = 859.05884
$wsCUL.Range("I132").Value = 625
$wsCUL.Range("J132").Value = 1067.1111
$wsCUL.Range("K132").Value = 5625
$wsCUL.Range("L132").Value = 9603.999900000001
$wsCUL.Range("M132").Value = -3095
$wsCUL.Range("N132").Value = -14663.9999

$wsGSM = $wb.Worksheets.Item("GSM")
# Row 80
$wsGSM.Range("H80").Value = 2791.6667
$wsGSM.Range("I80").Value = 3766.25
$wsGSM.Range("J80").Value = 2234.762
$wsGSM.Range("K80").Value = 3766.25
$wsGSM.Range("L80").Value = 2234.762
$wsGSM.Range("M80").Value = -2768.25
$wsGSM.Range("N80").Value = -4230.762000000001
# Row 83
$wsGSM.Range("H83").Value = 2791.6667
$wsGSM.Range("I83").Value = 3766.25
$wsGSM.Range("J83").Value = 2234.762
$wsGSM.Range("K83").Value = 18831.25
$wsGSM.Range("L83").Value = 11173.81
$wsGSM.Range("M83").Value = -13839.25
$wsGSM.Range("N83").Value = -21157.81
# Row 122
$wsGSM.Range("H122").Value = 1980.3334
$wsGSM.Range("I122").Value = 10007
$wsGSM.Range("J122").Value = 977
$wsGSM.Range("K122").Value = 30021
$wsGSM.Range("L122").Value = 2931
$wsGSM.Range("M122").Value = -27571
$wsGSM.Range("N122").Value = -7831
# Row 126
$wsGSM.Range("H126").Value = 1740.5814
$wsGSM.Range("I126").Value = 1452.4814
$wsGSM.Range("J126").Value = 2226.75
$wsGSM.Range("K126").Value = 4357.4442
$wsGSM.Range("L126").Value = 6680.25
$wsGSM.Range("M126").Value = -1887.4442
$wsGSM.Range("N126").Value = -11620.25

$wsLTW = $wb.Worksheets.Item("LTW")
# Row 132
$wsLTW.Range("H132").Value = 8934426
$wsLTW.Range("I132").Value = 16137518
$wsLTW.Range("J132").Value = 2592.28
$wsLTW.Range("K132").Value = 48412554
$wsLTW.Range("L132").Value = 7776.84
$wsLTW.Range("M132").Value = -48410024
$wsLTW.Range("N132").Value = -12836.84

$wsWVR = $wb.Worksheets.Item("WVR")
# Row 81
$wsWVR.Range("H81").Value = 35716120
$wsWVR.Range("I81").Value = 66668450
$wsWVR.Range("J81").Value = 1899.6923
$wsWVR.Range("K81").Value = 133336900
$wsWVR.Range("L81").Value = 3799.3846
$wsWVR.Range("M81").Value = -133335839
$wsWVR.Range("N81").Value = -5921.384599999999
# Row 84
$wsWVR.Range("H84").Value = 35716120
$wsWVR.Range("I84").Value = 66668450
$wsWVR.Range("J84").Value = 1899.6923
$wsWVR.Range("K84").Value = 666684500
$wsWVR.Range("L84").Value = 18996.923
$wsWVR.Range("M84").Value = -666679196
$wsWVR.Range("N84").Value = -29604.923
# Row 122
$wsWVR.Range("H122").Value = 1082.258
$wsWVR.Range("I122").Value = 907.875
$wsWVR.Range("J122").Value = 1268.2667
$wsWVR.Range("K122").Value = 2723.625
$wsWVR.Range("L122").Value = 3804.800099999999
$wsWVR.Range("M122").Value = -273.625
$wsWVR.Range("N122").Value = -8704.8001

